# Fix Training Data Issue (#48)
# The BF column ("Date") held a mangled literal ("5-6-2011-12") coming from
# the source file name instead of the actual game date. Correct it to the
# real ISO-ish date string "2012-05-06" for every data row (BF2:BF31).
# A leading apostrophe forces Excel to keep the value as plain text instead
# of re-interpreting "2012-05-06" as a date serial.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)  # column BF
    $cell.Value = "'2012-05-06"
}
